$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 (already has A1=1, B1=2)
$ws.Range("A1").Value = 1
$ws.Range("B1").Value = 2

# Row 2
$ws.Range("B2").Value = 1
$ws.Range("C2").Value = 2

# Row 5
$ws.Range("A5").Value = 1
$ws.Range("B5").Value = 2
$ws.Range("C5").Value = 3
$ws.Range("D5").Value = 4
$ws.Range("E5").Value = 5

# Update selection to match diff (activeCell F5)
$ws.Range("F5").Select()
